$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their textual formatting
# so numeric-looking strings like "1.000" or "5.050" are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.014.46'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.888.34'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '306.28'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("D7").Value = '0.5180'
$ws.Range("E7").Value = '  +2.83%  '
$ws.Range("D8").Value = '0.3755'
$ws.Range("E8").Value = '  +2.87%  '
$ws.Range("D9").Value = '0.07202'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").Value = '21.11'
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("D11").Value = '0.9011'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '0.07646'
$ws.Range("E12").Value = '  +2.00%  '
$ws.Range("D13").Value = '1.885.41'
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").Value = '94.26'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '5.242'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D17").Value = '0.000008496'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '27.062.74'
$ws.Range("D21").Value = '5.050'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").Value = '2.118.07'
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("D24").Value = '6.386'
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  +10.42%  '
$ws.Range("D26").Value = '146.17'
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '1.732'
$ws.Range("E27").Value = '  -2.74%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.06'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").Value = '114.21'
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").Value = '4.918'
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("D31").Value = '4.785'
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("D32").Value = '0.09202'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").Value = '0.05036'
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("E34").Value = '  +7.33%  '
$ws.Range("D35").Value = '0.7679'
$ws.Range("E35").Value = '  +2.35%  '
$ws.Range("D36").Value = '2.957'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '3.275'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").Value = '2.606'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").Value = '0.5589'
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").Value = '0.01985'
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").Value = '9.025'
$ws.Range("E42").Value = '  +5.18%  '
$ws.Range("D43").Value = '6.616'
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").Value = '118.80'
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("E45").Value = '  +2.16%  '
$ws.Range("D46").Value = '0.4828'
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.18'
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").Value = '1.601'
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("D50").Value = '37.67'
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("D51").Value = '63.98'
